# Applies the "added fault tolerance section to future work" edit.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "Given more time for this project, ..." paragraph: swap out
#    "things that we could use" for "aspects of our system that we
#    could improve".
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "things that we could use to enhance",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "aspects of our system that we could improve to enhance", 2)

# ---------------------------------------------------------------------
# 2. Split the "Currently, the both the serial ..." paragraph right
#    after "... would improve performance." and grow a brand new
#    paragraph describing fault tolerance in its place.  We build the
#    new paragraph's text by repeatedly locating the end of what was
#    just inserted (via a short, currently-unique anchor phrase) and
#    appending the next chunk after it - this keeps every insertion
#    anchored to real content instead of brittle character offsets.
# ---------------------------------------------------------------------

function Insert-After-Anchor($anchor, $text) {
    $rng = $d.Content
    $null = $rng.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $rng.Collapse(0) | Out-Null
    $rng.InsertAfter($text) | Out-Null
}

# 2a. Break the paragraph into two: start the new paragraph.
Insert-After-Anchor `
    "This would likely require additional memory but would improve performance." `
    "`rFault tolerance is another major aspect of distributed coWPAtty that could be improved.  Currently, the master node is a single point of failure.  By replicating the master node, system down time could be reduced since one of the master node replicas could be swapped in to act as the new master.  Additionally, another node could be added to act as the NFS host.  Relieving some of the pressure from the master node could possibly reduce the likelihood of failure.  "

# 2b. Continue building out the new paragraph's content.
Insert-After-Anchor `
    "Relieving some of the pressure from the master node could possibly reduce the likelihood of failure.  " `
    "Worker node failure is notably less critical than master node failure.  However, there is still room for improvement regarding the fault tolerance of the worker nodes.  Currently, the master node can detect worker node failure, but there is no automated system in place to deal with it.  A system administrator must restart the worker nodes to bring them all back up.  When the master detects worker node failure, it could send the restart command itself rather than having a system administrator perform that action.  Additionally, if the automated restart of the worker node by the master node failed, the master node could try to redistribute the work load to some of the worker nodes that were still alive.  Combining the automated restart of the worker nodes by the master node with the automated load redistribution would greatly enhance the fault tolerance in this system."

# ---------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark (Word's "last edit position" marker)
#    from its old spot (an empty paragraph before REFERENCES) to the
#    end of the brand new fault-tolerance paragraph, which is now the
#    location of the most recent edit.  Adding a bookmark that reuses
#    the reserved "_GoBack" name relocates the existing one instead of
#    creating a duplicate.
# ---------------------------------------------------------------------
$goBackAnchor = $d.Content
$null = $goBackAnchor.Find.Execute(
    "would greatly enhance the fault tolerance in this system.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$goBackAnchor.Collapse(0) | Out-Null
$d.Bookmarks.Add("_GoBack", $goBackAnchor) | Out-Null
